$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.797.87"
$ws.Range("E2").Value = "  +4.21%  "

$ws.Range("D3").Value = "1.867.32"
$ws.Range("E3").Value = "  +2.92%  "

$ws.Range("D4").Value = "'0.9994"
$ws.Range("E4").Value = "  -0.20%  "

$ws.Range("D5").Value = "'274.19"
$ws.Range("E5").Value = "  -1.42%  "

$ws.Range("E6").Value = "  -0.16%  "

$ws.Range("D7").Value = "'0.5299"
$ws.Range("E7").Value = "  +4.20%  "

$ws.Range("D8").Value = "'0.3389"
$ws.Range("E8").Value = "  -3.89%  "

$ws.Range("D9").Value = "'0.06813"
$ws.Range("E9").Value = "  +2.03%  "

$ws.Range("D10").Value = "'19.91"
$ws.Range("E10").Value = "  +0.15%  "

$ws.Range("D11").Value = "'0.7961"
$ws.Range("E11").Value = "  -2.94%  "

$ws.Range("D12").Value = "'0.07737"
$ws.Range("E12").Value = "  -1.54%  "

$ws.Range("D13").Value = "1.858.76"
$ws.Range("E13").Value = "  +1.87%  "

$ws.Range("D14").Value = "'90.14"

$ws.Range("D15").Value = "'5.129"
$ws.Range("E15").Value = "  +1.18%  "

$ws.Range("D16").Value = "'0.9988"
$ws.Range("E16").Value = "  -0.24%  "

$ws.Range("D17").Value = "'14.45"
$ws.Range("E17").Value = "  +2.60%  "

$ws.Range("D18").Value = "'0.000008013"
$ws.Range("E18").Value = "  -0.20%  "

$ws.Range("D19").Value = "'0.9991"
$ws.Range("E19").Value = "  -0.21%  "

$ws.Range("D20").Value = "26.822.64"
$ws.Range("E20").Value = "  +4.11%  "

$ws.Range("D21").Value = "2.099.11"
$ws.Range("E21").Value = "  +2.15%  "

$ws.Range("D22").Value = "'4.714"
$ws.Range("E22").Value = "  -0.60%  "

$ws.Range("D23").Value = "'9.973"
$ws.Range("E23").Value = "  -0.12%  "

$ws.Range("D24").Value = "'6.091"
$ws.Range("E24").Value = "  -0.31%  "

$ws.Range("D25").Value = "'2.374"
$ws.Range("E25").Value = "  +5.86%  "

$ws.Range("D26").Value = "'145.59"
$ws.Range("E26").Value = "  +2.27%  "

$ws.Range("B27").Value = "Toncoin"
$ws.Range("C27").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D27").Value = "'1.654"
$ws.Range("E27").Value = "  -0.79%  "

$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D28").Value = "'17.23"
$ws.Range("E28").Value = "  +0.70%  "

$ws.Range("D29").Value = "'112.61"
$ws.Range("E29").Value = "  +3.25%  "

$ws.Range("E30").Value = "  +0.06%  "

$ws.Range("D31").Value = "'4.318"
$ws.Range("E31").Value = "  +2.41%  "

$ws.Range("D32").Value = "'0.08870"
$ws.Range("E32").Value = "  +1.55%  "

$ws.Range("D33").Value = "'0.04926"
$ws.Range("E33").Value = "  +1.36%  "

$ws.Range("D34").Value = "'1.165"
$ws.Range("E34").Value = "  +3.20%  "

$ws.Range("D35").Value = "'0.7282"
$ws.Range("E35").Value = "  +0.15%  "

$ws.Range("D36").Value = "'2.875"
$ws.Range("E36").Value = "  -0.64%  "

$ws.Range("D37").Value = "'3.211"
$ws.Range("E37").Value = "  +1.94%  "

$ws.Range("D38").Value = "'2.335"
$ws.Range("E38").Value = "  -1.21%  "

$ws.Range("D39").Value = "'0.01851"
$ws.Range("E39").Value = "  +0.20%  "

$ws.Range("D40").Value = "'0.5097"
$ws.Range("E40").Value = "  -0.67%  "

$ws.Range("B41").Value = "Quant"
$ws.Range("C41").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D41").Value = "'116.44"
$ws.Range("E41").Value = "  +2.02%  "

$ws.Range("B42").Value = "TrustWalletToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D42").Value = "'0.9359"
$ws.Range("E42").Value = "  -3.42%  "

$ws.Range("D43").Value = "'6.121"
$ws.Range("E43").Value = "  -1.57%  "

$ws.Range("D44").Value = "'8.015"
$ws.Range("E44").Value = "  +0.08%  "

$ws.Range("D45").Value = "'0.9988"
$ws.Range("E45").Value = "  -0.21%  "

$ws.Range("E46").Value = "  -2.37%  "

$ws.Range("E47").Value = "  -2.75%  "

$ws.Range("D48").Value = "'9.271"
$ws.Range("E48").Value = "  +0.61%  "

$ws.Range("E49").Value = "  -0.81%  "

$ws.Range("D50").Value = "'0.05947"
$ws.Range("E50").Value = "  +2.19%  "

$ws.Range("E51").Value = "  -1.74%  "
